$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# D-column values are numeric-looking text (e.g. "534.31"); force them to
# remain plain text (matching the original inline-string cells) by applying
# a text number format before assignment, then clearing the format so the
# cell keeps no explicit style, just like the source file.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.936.01"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.128.53"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.31"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.21"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.485"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +7.26%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.413"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.07%  "
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.667.70"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.98"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.024.65"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.126.03"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.85"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.15"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.69"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +7.56%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.71"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.55"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.54"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.17"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.50%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.56"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.17"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("E37").Value = "  +4.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.76"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0674"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.570.61"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.13%  "
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.702"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "37.82"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.15"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.83"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.750"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0926"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.57%  "
